# New order came in (Ajay Dwarkunde, 2026-01-13 10:51) - insert it at the top of the
# "All Orders" log (rows are newest-first) and bump the "Daily Summary" order count.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("All Orders")

# Push the existing order rows down one slot to make room for the new order at row 2.
$ws.Rows.Item(2).Insert()

$ws.Range("A2").Value2 = 4
$ws.Range("B2").Value2 = "2026-01-13 10:51"
$ws.Range("C2").Value2 = "Ajay Dwarkunde"
# Phone numbers are kept as text in this sheet (leading apostrophe keeps it text, not numeric).
$ws.Range("D2").Value2 = "'8087172173"
$ws.Range("E2").Value2 = "wakad, pune 411057"
$ws.Range("F2").Value2 = "Girl Holding Hands Thali x1, Kalash Haldi Kunku (Golden) x1, Kite Haldi Kunku Set x1"
$ws.Range("G2").Value2 = 0
$ws.Range("H2").Value2 = "NEW"
$ws.Range("I2").Value2 = "PENDING"

# Bump the running order total on the Daily Summary tab.
$summary = $wb.Worksheets.Item("Daily Summary")
$summary.Range("B2").Value2 = 4
